# Reorder data rows (2-93) according to the source mapping derived from the diff.
# Entire row content (across all used columns) moves as a unit, keyed by the record's
# "Id" value in column A. Startdatum/Slutdatum (columns Y/AA) are excluded from the
# rewrite because every row shares the same literal text value "2026-02-15"; touching
# those cells via Value2 would coerce the stored text into a date serial number, which
# is not part of the intended change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = @{}
$srcRow[2] = 4
$srcRow[3] = 3
$srcRow[4] = 5
$srcRow[5] = 2
$srcRow[6] = 6
$srcRow[7] = 8
$srcRow[8] = 7
$srcRow[9] = 9
$srcRow[10] = 10
$srcRow[11] = 12
$srcRow[12] = 13
$srcRow[13] = 14
$srcRow[14] = 11
$srcRow[15] = 17
$srcRow[16] = 15
$srcRow[17] = 23
$srcRow[18] = 20
$srcRow[19] = 21
$srcRow[20] = 16
$srcRow[21] = 19
$srcRow[22] = 22
$srcRow[23] = 18
$srcRow[24] = 26
$srcRow[25] = 24
$srcRow[26] = 29
$srcRow[27] = 25
$srcRow[28] = 28
$srcRow[29] = 27
$srcRow[30] = 30
$srcRow[31] = 31
$srcRow[32] = 33
$srcRow[33] = 34
$srcRow[34] = 35
$srcRow[35] = 36
$srcRow[36] = 32
$srcRow[37] = 39
$srcRow[38] = 37
$srcRow[39] = 40
$srcRow[40] = 38
$srcRow[41] = 42
$srcRow[42] = 43
$srcRow[43] = 41
$srcRow[44] = 44
$srcRow[45] = 45
$srcRow[46] = 50
$srcRow[47] = 52
$srcRow[48] = 49
$srcRow[49] = 47
$srcRow[50] = 48
$srcRow[51] = 46
$srcRow[52] = 51
$srcRow[53] = 54
$srcRow[54] = 55
$srcRow[55] = 56
$srcRow[56] = 57
$srcRow[57] = 53
$srcRow[58] = 62
$srcRow[59] = 59
$srcRow[60] = 58
$srcRow[61] = 63
$srcRow[62] = 61
$srcRow[63] = 60
$srcRow[64] = 64
$srcRow[65] = 65
$srcRow[66] = 66
$srcRow[67] = 69
$srcRow[68] = 70
$srcRow[69] = 68
$srcRow[70] = 67
$srcRow[71] = 78
$srcRow[72] = 80
$srcRow[73] = 71
$srcRow[74] = 82
$srcRow[75] = 81
$srcRow[76] = 76
$srcRow[77] = 75
$srcRow[78] = 72
$srcRow[79] = 73
$srcRow[80] = 79
$srcRow[81] = 74
$srcRow[82] = 77
$srcRow[83] = 84
$srcRow[84] = 85
$srcRow[85] = 83
$srcRow[86] = 89
$srcRow[87] = 86
$srcRow[88] = 87
$srcRow[89] = 88
$srcRow[90] = 90
$srcRow[91] = 92
$srcRow[92] = 93
$srcRow[93] = 91

$firstRow = 2
$lastRow = 93

# Column blocks to rewrite, skipping Y (25) and AA (27):
#   A:X  -> columns 1..24
#   Z    -> column 26
#   AB:AY-> columns 28..51
$blocks = @(
  @{ StartCol = 1;  EndCol = 24 },
  @{ StartCol = 26; EndCol = 26 },
  @{ StartCol = 28; EndCol = 51 }
)

foreach ($block in $blocks) {
  $startCol = $block.StartCol
  $endCol = $block.EndCol
  $numCols = $endCol - $startCol + 1
  $numRows = $lastRow - $firstRow + 1

  $srcRange = $ws.Range($ws.Cells.Item($firstRow, $startCol), $ws.Cells.Item($lastRow, $endCol))
  $data = $srcRange.Value2

  $new = New-Object 'object[,]' $numRows,$numCols

  for ($r = 0; $r -lt $numRows; $r++) {
    $destRowNum = $firstRow + $r
    $sourceRowNum = $srcRow[$destRowNum]
    $sourceIdx = $sourceRowNum - $firstRow + 1
    for ($c = 1; $c -le $numCols; $c++) {
      $new[$r, $c - 1] = $data[$sourceIdx, $c]
    }
  }

  $destRange = $ws.Range($ws.Cells.Item($firstRow, $startCol), $ws.Cells.Item($lastRow, $endCol))
  $destRange.Value2 = $new
}
